$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly added state data (Louisiana, Maine, and Ohio through Wyoming) ---

# Row 19: LOUISIANA
$ws.Range("B19").Value = 57
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = "<3"
$ws.Range("J19").Value = 8
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 31
$ws.Range("M19").Value = 21
$ws.Range("N19").Value = 51
$ws.Range("O19").Value = "<3"
$ws.Range("P19").Value = "<3"

# Row 20: MAINE
$ws.Range("B20").Value = 7
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = "<3"
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = "<3"
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = "<3"
$ws.Range("N20").Value = 6
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = "<3"

# Row 36: OHIO
$ws.Range("B36").Value = 233
$ws.Range("C36").Value = 18
$ws.Range("D36").Value = 8
$ws.Range("E36").Value = 30
$ws.Range("F36").Value = 289
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 4
$ws.Range("I36").Value = 12
$ws.Range("J36").Value = 6
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 172
$ws.Range("M36").Value = 53
$ws.Range("N36").Value = 218
$ws.Range("O36").Value = 8
$ws.Range("P36").Value = "<3"

# Row 37: OKLAHOMA
$ws.Range("B37").Value = 28
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 4
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 42
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = "<3"
$ws.Range("I37").Value = 5
$ws.Range("J37").Value = 6
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 18
$ws.Range("M37").Value = 10
$ws.Range("N37").Value = 27
$ws.Range("O37").Value = "<3"
$ws.Range("P37").Value = 0

# Row 38: OREGON
$ws.Range("B38").Value = 44
$ws.Range("C38").Value = 8
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 3
$ws.Range("F38").Value = 56
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = "<3"
$ws.Range("I38").Value = 7
$ws.Range("J38").Value = "<3"
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 32
$ws.Range("M38").Value = 11
$ws.Range("N38").Value = 42
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = "<3"

# Row 39: PENNSYLVANIA
$ws.Range("B39").Value = 80
$ws.Range("C39").Value = 15
$ws.Range("D39").Value = 6
$ws.Range("E39").Value = 5
$ws.Range("F39").Value = 106
$ws.Range("G39").Value = 14
$ws.Range("H39").Value = "<3"
$ws.Range("I39").Value = 9
$ws.Range("J39").Value = 7
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 53
$ws.Range("M39").Value = 20
$ws.Range("N39").Value = 68
$ws.Range("O39").Value = 4
$ws.Range("P39").Value = "<3"

# Row 40: RHODE ISLAND
$ws.Range("B40").Value = 5
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 8
$ws.Range("G40").Value = "<3"
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = "<3"
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5
$ws.Range("M40").Value = "<3"
$ws.Range("N40").Value = 5
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0

# Row 41: SOUTH CAROLINA
$ws.Range("B41").Value = 42
$ws.Range("C41").Value = 14
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 61
$ws.Range("G41").Value = 13
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 6
$ws.Range("J41").Value = 11
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 26
$ws.Range("M41").Value = 22
$ws.Range("N41").Value = 39
$ws.Range("O41").Value = 3
$ws.Range("P41").Value = 0

# Row 42: SOUTH DAKOTA
$ws.Range("B42").Value = 9
$ws.Range("C42").Value = 5
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 15
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 4
$ws.Range("J42").Value = "<3"
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 6
$ws.Range("M42").Value = 3
$ws.Range("N42").Value = 9
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0

# Row 43: TENNESSEE
$ws.Range("B43").Value = 53
$ws.Range("C43").Value = 12
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 3
$ws.Range("F43").Value = 69
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = "<3"
$ws.Range("I43").Value = 6
$ws.Range("J43").Value = 7
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28
$ws.Range("M43").Value = 23
$ws.Range("N43").Value = 51
$ws.Range("O43").Value = 3
$ws.Range("P43").Value = "<3"

# Row 44: TEXAS
$ws.Range("B44").Value = 337
$ws.Range("C44").Value = 63
$ws.Range("D44").Value = 16
$ws.Range("E44").Value = 17
$ws.Range("F44").Value = 433
$ws.Range("G44").Value = 50
$ws.Range("H44").Value = 16
$ws.Range("I44").Value = 32
$ws.Range("J44").Value = 34
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 205
$ws.Range("M44").Value = 123
$ws.Range("N44").Value = 312
$ws.Range("O44").Value = 11
$ws.Range("P44").Value = "<3"

# Row 45: UTAH
$ws.Range("B45").Value = 15
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 4
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 22
$ws.Range("G45").Value = "<3"
$ws.Range("H45").Value = "<3"
$ws.Range("I45").Value = 3
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 10
$ws.Range("M45").Value = 4
$ws.Range("N45").Value = 14
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0

# Row 46: VERMONT
$ws.Range("B46").Value = 3
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 5
$ws.Range("G46").Value = "<3"
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = "<3"
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = "<3"
$ws.Range("M46").Value = "<3"
$ws.Range("N46").Value = 3
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0

# Row 47: VIRGINIA
$ws.Range("B47").Value = 7
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 10
$ws.Range("G47").Value = "<3"
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = "<3"
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 3
$ws.Range("M47").Value = 3
$ws.Range("N47").Value = 7
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0

# Row 48: WASHINGTON
$ws.Range("B48").Value = 86
$ws.Range("C48").Value = 13
$ws.Range("D48").Value = 7
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = 111
$ws.Range("G48").Value = 11
$ws.Range("H48").Value = "<3"
$ws.Range("I48").Value = 8
$ws.Range("J48").Value = 4
$ws.Range("K48").Value = "<3"
$ws.Range("L48").Value = 48
$ws.Range("M48").Value = 26
$ws.Range("N48").Value = 81
$ws.Range("O48").Value = 5
$ws.Range("P48").Value = 0

# Row 49: WEST VIRGINIA
$ws.Range("B49").Value = 7
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 1
$ws.Range("F49").Value = 10
$ws.Range("G49").Value = "<3"
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = "<3"
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3
$ws.Range("M49").Value = 3
$ws.Range("N49").Value = 7
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0

# Row 50: WISCONSIN
$ws.Range("B50").Value = 45
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 50
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = "<3"
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 35
$ws.Range("M50").Value = 11
$ws.Range("N50").Value = 43
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0

# Row 51: WYOMING
$ws.Range("B51").Value = 1
$ws.Range("C51").Value = 4
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 1
$ws.Range("F51").Value = 6
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 3
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = "<3"
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = "<3"
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0

# --- Row 35 (NORTH DAKOTA) previously had a formatted-but-empty cell in C35; clear it fully ---
$ws.Range("C35").Clear()

# --- Update the active selection / last worked cell to reflect where data entry left off ---
$ws.Range("A36").Select()
